$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 422 (shifts existing rows 422:486 down to 423:487)
$ws.Rows("422").Insert()

# Populate the newly inserted row with the new price observation
$ws.Range("A422").Value = 10
$ws.Range("B422").Value = "Vega Modelo de Temuco"
$ws.Range("C422").Value = "La Araucanía"
$ws.Range("D422").Value = 44491
$ws.Range("E422").Value = 9
$ws.Range("F422").Value = 100112021
$ws.Range("G422").Value = "Ají"
$ws.Range("H422").Value = "Inferno"
$ws.Range("I422").Value = "Primera"
$ws.Range("J422").Value = 75
$ws.Range("K422").Value = 35000
$ws.Range("L422").Value = 35000
$ws.Range("M422").Value = 35000
$ws.Range("N422").Value = "$/caja 15 kilos"
$ws.Range("O422").Value = "Región de Arica y Parinacota"
$ws.Range("P422").Value = 2333
$ws.Range("Q422").Value = 15
$ws.Range("R422").Value = "Hortaliza"
